$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translations")

# Style template cells (existing, unmodified rows) used purely as formatting sources.
# Row 515 col C uses the 'wrap text' style; row 45 col B uses the 'plain / no-wrap' style.
# Both are far away from the edited range (516-535) and are left untouched by this script.
$wrapStyleSrc = $ws.Cells.Item(515, 3)
$plainStyleSrc = $ws.Cells.Item(45, 2)

# Row 516
$wrapStyleSrc.Copy($ws.Cells.Item(516, 1))
$ws.Cells.Item(516, 1).ClearContents()
$wrapStyleSrc.Copy($ws.Cells.Item(516, 2))
$ws.Cells.Item(516, 2).Value = "common._each_session"
$wrapStyleSrc.Copy($ws.Cells.Item(516, 3))
$ws.Cells.Item(516, 3).Value = "*Jede Sitzung schliesst mit einer Challenge ab, betrachtet Sie als eine Art Hausaufgabe oder als ein neues Werkzeug, dass ihr in Zukunft benutzen könnt."
$wrapStyleSrc.Copy($ws.Cells.Item(516, 4))
$ws.Cells.Item(516, 4).Value = "*Each session concludes with a challenge, think of it as a homework assignment or a new tool you can use in the future."
$wrapStyleSrc.Copy($ws.Cells.Item(516, 5))
$ws.Cells.Item(516, 5).Value = "*Cada sessão termina com um desafio, pense nele como uma tarefa de casa ou uma nova ferramenta que você poderá usar no futuro."

# Row 517
$plainStyleSrc.Copy($ws.Cells.Item(517, 2))
$ws.Cells.Item(517, 2).Value = "sessions.questions_in_sessions_are_on"
$plainStyleSrc.Copy($ws.Cells.Item(517, 3))
$ws.Cells.Item(517, 3).Value = "Der Fokus der Fragen in dieser Sitzung liegt auf"
$plainStyleSrc.Copy($ws.Cells.Item(517, 4))
$ws.Cells.Item(517, 4).Value = "The focus of the questions in this sessions are on"
$plainStyleSrc.Copy($ws.Cells.Item(517, 5))
$ws.Cells.Item(517, 5).Value = "O foco das questões nesta sessão está em"

# Row 518
$plainStyleSrc.Copy($ws.Cells.Item(518, 2))
$ws.Cells.Item(518, 2).Value = "sessions.a_handpicked_challenge_is_part_1"
$plainStyleSrc.Copy($ws.Cells.Item(518, 3))
$ws.Cells.Item(518, 3).Value = "Eine Auswahl an ausgewählten Core - Challenges der Kategorie Übung"
$plainStyleSrc.Copy($ws.Cells.Item(518, 4))
$ws.Cells.Item(518, 4).Value = "A set of handpicked Core - Challenges of the category Exercise"
$plainStyleSrc.Copy($ws.Cells.Item(518, 5))
$ws.Cells.Item(518, 5).Value = "Um conjunto de Desafios Core cuidadosamente selecionados da categoria Exercício"

# Row 519
$plainStyleSrc.Copy($ws.Cells.Item(519, 2))
$ws.Cells.Item(519, 2).Value = "sessions.a_handpicked_challenge_is_part_2"
$plainStyleSrc.Copy($ws.Cells.Item(519, 3))
$ws.Cells.Item(519, 3).Value = "rundet das Bild ab."
$plainStyleSrc.Copy($ws.Cells.Item(519, 4))
$ws.Cells.Item(519, 4).Value = "is rounding up the picture."
$plainStyleSrc.Copy($ws.Cells.Item(519, 5))
$ws.Cells.Item(519, 5).Value = "está completando o quadro."

# Row 520
$plainStyleSrc.Copy($ws.Cells.Item(520, 2))
$ws.Cells.Item(520, 2).Value = "sessions.picked_special_challenge_part_1"
$plainStyleSrc.Copy($ws.Cells.Item(520, 3))
$ws.Cells.Item(520, 3).Value = "Die besonderen Herausforderungen in dieser Session:"
$plainStyleSrc.Copy($ws.Cells.Item(520, 4))
$ws.Cells.Item(520, 4).Value = "The Special Challenges in this session:"
$plainStyleSrc.Copy($ws.Cells.Item(520, 5))
$ws.Cells.Item(520, 5).Value = "Os desafios especiais nesta sessão:"

# Row 521
$plainStyleSrc.Copy($ws.Cells.Item(521, 2))
$ws.Cells.Item(521, 2).Value = "therapists.programs_offered_by"
$plainStyleSrc.Copy($ws.Cells.Item(521, 3))
$ws.Cells.Item(521, 3).Value = "Programme angeboten von"
$plainStyleSrc.Copy($ws.Cells.Item(521, 4))
$ws.Cells.Item(521, 4).Value = "Programs offered by"
$plainStyleSrc.Copy($ws.Cells.Item(521, 5))
$ws.Cells.Item(521, 5).Value = "Programas oferecidos por"

# Row 522
$plainStyleSrc.Copy($ws.Cells.Item(522, 2))
$ws.Cells.Item(522, 2).Value = "therapists.contacts"
$plainStyleSrc.Copy($ws.Cells.Item(522, 3))
$ws.Cells.Item(522, 3).Value = "Kontakte"
$plainStyleSrc.Copy($ws.Cells.Item(522, 4))
$ws.Cells.Item(522, 4).Value = "Contacts"
$plainStyleSrc.Copy($ws.Cells.Item(522, 5))
$ws.Cells.Item(522, 5).Value = "Contatos"

# Row 523
$plainStyleSrc.Copy($ws.Cells.Item(523, 2))
$ws.Cells.Item(523, 2).Value = "therapist.contact_linkedin"
$plainStyleSrc.Copy($ws.Cells.Item(523, 3))
$ws.Cells.Item(523, 3).Value = "LinkedIn"
$plainStyleSrc.Copy($ws.Cells.Item(523, 4))
$ws.Cells.Item(523, 4).Value = "LinkedIn"
$plainStyleSrc.Copy($ws.Cells.Item(523, 5))
$ws.Cells.Item(523, 5).Value = "LinkedIn"

# Row 524
$plainStyleSrc.Copy($ws.Cells.Item(524, 2))
$ws.Cells.Item(524, 2).Value = "therapist.contact_phone"
$plainStyleSrc.Copy($ws.Cells.Item(524, 3))
$ws.Cells.Item(524, 3).Value = "Telefon"
$plainStyleSrc.Copy($ws.Cells.Item(524, 4))
$ws.Cells.Item(524, 4).Value = "Phone"
$plainStyleSrc.Copy($ws.Cells.Item(524, 5))
$ws.Cells.Item(524, 5).Value = "Telefone"

# Row 525
$plainStyleSrc.Copy($ws.Cells.Item(525, 2))
$ws.Cells.Item(525, 2).Value = "therapist.contact_website"
$plainStyleSrc.Copy($ws.Cells.Item(525, 3))
$ws.Cells.Item(525, 3).Value = "Webseite"
$plainStyleSrc.Copy($ws.Cells.Item(525, 4))
$ws.Cells.Item(525, 4).Value = "Website"
$plainStyleSrc.Copy($ws.Cells.Item(525, 5))
$ws.Cells.Item(525, 5).Value = "Site"

# Row 526
$plainStyleSrc.Copy($ws.Cells.Item(526, 2))
$ws.Cells.Item(526, 2).Value = "therapist.contact_email"
$plainStyleSrc.Copy($ws.Cells.Item(526, 3))
$ws.Cells.Item(526, 3).Value = "Email"
$plainStyleSrc.Copy($ws.Cells.Item(526, 4))
$ws.Cells.Item(526, 4).Value = "Email"
$plainStyleSrc.Copy($ws.Cells.Item(526, 5))
$ws.Cells.Item(526, 5).Value = "Email"

# Row 527
$plainStyleSrc.Copy($ws.Cells.Item(527, 2))
$ws.Cells.Item(527, 2).Value = "therapist.contact_instagram"
$plainStyleSrc.Copy($ws.Cells.Item(527, 3))
$ws.Cells.Item(527, 3).Value = "Instagram"
$plainStyleSrc.Copy($ws.Cells.Item(527, 4))
$ws.Cells.Item(527, 4).Value = "Instagram"
$plainStyleSrc.Copy($ws.Cells.Item(527, 5))
$ws.Cells.Item(527, 5).Value = "Instagram"

# Row 528
$plainStyleSrc.Copy($ws.Cells.Item(528, 2))
$ws.Cells.Item(528, 2).Value = "therapist.contact_youtube"
$plainStyleSrc.Copy($ws.Cells.Item(528, 3))
$ws.Cells.Item(528, 3).Value = "Youtube"
$plainStyleSrc.Copy($ws.Cells.Item(528, 4))
$ws.Cells.Item(528, 4).Value = "Youtube"
$plainStyleSrc.Copy($ws.Cells.Item(528, 5))
$ws.Cells.Item(528, 5).Value = "Youtube"

# Row 529
$plainStyleSrc.Copy($ws.Cells.Item(529, 2))
$ws.Cells.Item(529, 2).Value = "languages.en"
$plainStyleSrc.Copy($ws.Cells.Item(529, 3))
$ws.Cells.Item(529, 3).Value = "Englisch"
$plainStyleSrc.Copy($ws.Cells.Item(529, 4))
$ws.Cells.Item(529, 4).Value = "English"
$plainStyleSrc.Copy($ws.Cells.Item(529, 5))
$ws.Cells.Item(529, 5).Value = "Inglês"

# Row 530
$plainStyleSrc.Copy($ws.Cells.Item(530, 2))
$ws.Cells.Item(530, 2).Value = "languages.de"
$ws.Cells.Item(530, 2).NumberFormat = "@"
$plainStyleSrc.Copy($ws.Cells.Item(530, 3))
$ws.Cells.Item(530, 3).Value = "Deutsch"
$plainStyleSrc.Copy($ws.Cells.Item(530, 4))
$ws.Cells.Item(530, 4).Value = "German"
$plainStyleSrc.Copy($ws.Cells.Item(530, 5))
$ws.Cells.Item(530, 5).Value = "Alemão"

# Row 531
$plainStyleSrc.Copy($ws.Cells.Item(531, 2))
$ws.Cells.Item(531, 2).Value = "languages.pt"
$ws.Cells.Item(531, 2).NumberFormat = "@"
$plainStyleSrc.Copy($ws.Cells.Item(531, 3))
$ws.Cells.Item(531, 3).Value = "Portugiesisch"
$plainStyleSrc.Copy($ws.Cells.Item(531, 4))
$ws.Cells.Item(531, 4).Value = "Portuguese"
$plainStyleSrc.Copy($ws.Cells.Item(531, 5))
$ws.Cells.Item(531, 5).Value = "Português"

# Row 532
$plainStyleSrc.Copy($ws.Cells.Item(532, 2))
$ws.Cells.Item(532, 2).Value = "languages.es"
$ws.Cells.Item(532, 2).NumberFormat = "@"
$plainStyleSrc.Copy($ws.Cells.Item(532, 3))
$ws.Cells.Item(532, 3).Value = "Spanisch"
$plainStyleSrc.Copy($ws.Cells.Item(532, 4))
$ws.Cells.Item(532, 4).Value = "Spanish"
$plainStyleSrc.Copy($ws.Cells.Item(532, 5))
$ws.Cells.Item(532, 5).Value = "Espanhol"

# Row 533
$plainStyleSrc.Copy($ws.Cells.Item(533, 2))
$ws.Cells.Item(533, 2).Value = "languages.fr"
$ws.Cells.Item(533, 2).NumberFormat = "@"
$plainStyleSrc.Copy($ws.Cells.Item(533, 3))
$ws.Cells.Item(533, 3).Value = "Französisch"
$plainStyleSrc.Copy($ws.Cells.Item(533, 4))
$ws.Cells.Item(533, 4).Value = "French"
$plainStyleSrc.Copy($ws.Cells.Item(533, 5))
$ws.Cells.Item(533, 5).Value = "Francês"

# Row 534
$plainStyleSrc.Copy($ws.Cells.Item(534, 2))
$ws.Cells.Item(534, 2).Value = "languages.sv"
$ws.Cells.Item(534, 2).NumberFormat = "@"
$plainStyleSrc.Copy($ws.Cells.Item(534, 3))
$ws.Cells.Item(534, 3).Value = "Schwedisch"
$plainStyleSrc.Copy($ws.Cells.Item(534, 4))
$ws.Cells.Item(534, 4).Value = "Swedish"
$plainStyleSrc.Copy($ws.Cells.Item(534, 5))
$ws.Cells.Item(534, 5).Value = "Sueco"

# Row 535
$plainStyleSrc.Copy($ws.Cells.Item(535, 2))
$ws.Cells.Item(535, 2).Value = "languages.gr"
$ws.Cells.Item(535, 2).NumberFormat = "@"
$plainStyleSrc.Copy($ws.Cells.Item(535, 3))
$ws.Cells.Item(535, 3).Value = "Griechisch"
$plainStyleSrc.Copy($ws.Cells.Item(535, 4))
$ws.Cells.Item(535, 4).Value = "Greek"
$plainStyleSrc.Copy($ws.Cells.Item(535, 5))
$ws.Cells.Item(535, 5).Value = "Grego"
